# Insert a new data row at row 41 (pushes existing rows 41-122 down to 42-123),
# matching the commit "Fruta / hortaliza, semanal" which adds a new weekly
# observation for Arveja Verde at Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 41, shifting rows 41:122 down to 42:123.
# Excel copies formatting (including the date style on column D) from the
# row above automatically.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new observation values.
$ws.Range("A41").Value2 = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value2 = 45259
$ws.Range("E41").Value2 = 16
$ws.Range("F41").Value2 = 100112022
$ws.Range("G41").Value = "Arveja Verde"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value2 = 200
$ws.Range("K41").Value2 = 25000
$ws.Range("L41").Value2 = 25000
$ws.Range("M41").Value2 = 25000
$ws.Range("N41").Value = "$/saco 25 kilos"
$ws.Range("O41").Value = "Región del Maule"
$ws.Range("P41").Value2 = 1000
$ws.Range("Q41").Value2 = 25
$ws.Range("R41").Value = "Hortaliza"
